$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,11
$arr[0,0] = 0.1233546200258786
$arr[0,1] = 0.07308738474181808
$arr[0,2] = 0.1174230245179544
$arr[0,3] = 2.628466502658767
$arr[0,4] = 1.935490001332639
$arr[0,5] = 1.675732621926272
$arr[0,6] = 1.88147098425808
$arr[0,7] = 0.1896188743661629
$arr[0,8] = 1.91794182479407
$arr[0,9] = 0.1847726683820952
$arr[0,10] = 0.5634616086912274
$arr[1,0] = 0.1216966938540196
$arr[1,1] = 0.07157544687918005
$arr[1,2] = 0.1175862777724745
$arr[1,3] = 2.646582355138428
$arr[1,4] = 1.95137172918615
$arr[1,5] = 1.690331141138813
$arr[1,6] = 1.895018260411909
$arr[1,7] = 0.1910742614136449
$arr[1,8] = 1.7831272642523
$arr[1,9] = 0.1858199146600032
$arr[1,10] = 0.5389929666983235
$arr[2,0] = 0.1207101151223995
$arr[2,1] = 0.07065570940591215
$arr[2,2] = 0.1177191523100607
$arr[2,3] = 2.659192352163402
$arr[2,4] = 1.962401979382591
$arr[2,5] = 1.700135606189932
$arr[2,6] = 1.904416877475327
$arr[2,7] = 0.1920321914351142
$arr[2,8] = 1.70070984224381
$arr[2,9] = 0.1865139022580777
$arr[2,10] = 0.5241218064352751
$arr[3,0] = 0.1203160170097988
$arr[3,1] = 0.07028310570553487
$arr[3,2] = 0.1177815265374829
$arr[3,3] = 2.664704628336111
$arr[3,4] = 1.967217975753826
$arr[3,5] = 1.704342384382215
$arr[3,6] = 1.908518381984841
$arr[3,7] = 0.1924387432598458
$arr[3,8] = 1.667216053671865
$arr[3,9] = 0.1868095464587221
$arr[3,10] = 0.5181004295693725
$arr[4,0] = 0.1202510582928653
$arr[4,1] = 0.07022136885501595
$arr[4,2] = 0.1177923812004718
$arr[4,3] = 2.665642494378787
$arr[4,4] = 1.96803704649642
$arr[4,5] = 1.705053681273057
$arr[4,6] = 1.909215824774385
$arr[4,7] = 0.1925072290018424
$arr[4,8] = 1.661660039888517
$arr[4,9] = 0.1868594139854487
$arr[4,10] = 0.5171029347883334
$arr[5,0] = 0.1207047679616622
$arr[5,1] = 0.07065067540018077
$arr[5,2] = 0.1177199601738153
$arr[5,3] = 2.659265180316815
$arr[5,4] = 1.962465630206196
$arr[5,5] = 1.700191484514136
$arr[5,6] = 1.904471092834925
$arr[5,7] = 0.1920376087733064
$arr[5,8] = 1.700257758678561
$arr[5,9] = 0.1865178374087115
$arr[5,10] = 0.5240404427124048
$arr[6,0] = 0.1227764780463403
$arr[6,1] = 0.07256430403436553
$arr[6,2] = 0.1174725524421678
$arr[6,3] = 2.634404103238296
$arr[6,4] = 1.940700311402949
$arr[6,5] = 1.680591613505612
$arr[6,6] = 1.885917687752915
$arr[6,7] = 0.1901073517873346
$arr[6,8] = 1.871384131065213
$arr[6,9] = 0.1851231916746716
$arr[6,10] = 0.5549933298161136
$arr[7,0] = 0.1270864645828595
$arr[7,1] = 0.07638385856568419
$arr[7,2] = 0.1172455207832126
$arr[7,3] = 2.597461877960569
$arr[7,4] = 1.90818808463267
$arr[7,5] = 1.648832296745354
$arr[7,6] = 1.858118747337386
$arr[7,7] = 0.1868317816121383
$arr[7,8] = 2.209759800880249
$arr[7,9] = 0.1827918357039557
$arr[7,10] = 0.6168916893543823
$arr[8,0] = 0.1304019288635345
$arr[8,1] = 0.07922949568546045
$arr[8,2] = 0.1172350739304591
$arr[8,3] = 2.577540505782707
$arr[8,4] = 1.890533656297478
$arr[8,5] = 1.629574073958878
$arr[8,6] = 1.842944077109905
$arr[8,7] = 0.1847349980985662
$arr[8,8] = 2.460026056748632
$arr[8,9] = 0.1813237572961981
$arr[8,10] = 0.6630895905744154
$arr[9,0] = 0.1319421923508344
$arr[9,1] = 0.08053234179728008
$arr[9,2] = 0.1172640661775315
$arr[9,3] = 2.570050070415689
$arr[9,4] = 1.883862646612613
$arr[9,5] = 1.621699193849565
$arr[9,6] = 1.837184124079677
$arr[9,7] = 0.1838481890821484
$arr[9,8] = 2.574232402674568
$arr[9,9] = 0.1807087796288158
$arr[9,10] = 0.6842610552834714
$arr[10,0] = 0.1325300193726946
$arr[10,1] = 0.08102686855951902
$arr[10,2] = 0.1172798789924911
$arr[10,3] = 2.56744002946995
$arr[10,4] = 1.881532660998076
$arr[10,5] = 1.618844668469208
$arr[10,6] = 1.8351676362633
$arr[10,7] = 0.1835220024005828
$arr[10,8] = 2.617529852191694
$arr[10,9] = 0.1804834844335055
$arr[10,10] = 0.6923002916504828
$arr[11,0] = 0.1324032180896637
$arr[11,1] = 0.08092031199195304
$arr[11,2] = 0.1172762587131384
$arr[11,3] = 2.567992072006831
$arr[11,4] = 1.88202573026247
$arr[11,5] = 1.619453767710297
$arr[11,6] = 1.83559459382559
$arr[11,7] = 0.1835918244587447
$arr[11,8] = 2.608202774351923
$arr[11,9] = 0.1805316687836438
$arr[11,10] = 0.6905679224322938
$arr[12,0] = 0.1319904620453229
$arr[12,1] = 0.08057300360984243
$arr[12,2] = 0.117265270329618
$arr[12,3] = 2.569830800668058
$arr[12,4] = 1.883667022074107
$arr[12,5] = 1.621461793777158
$arr[12,6] = 1.837014924185027
$arr[12,7] = 0.183821160591485
$arr[12,8] = 2.577793517385487
$arr[12,9] = 0.1806900925728918
$arr[12,10] = 0.6849220078813829
$arr[13,0] = 0.1317382299205576
$arr[13,1] = 0.08036041821362261
$arr[13,2] = 0.1172591686257256
$arr[13,3] = 2.570986573737542
$arr[13,4] = 1.884697927184334
$arr[13,5] = 1.622708379070971
$arr[13,6] = 1.837906373079797
$arr[13,7] = 0.1839628892367422
$arr[13,8] = 2.559173430970304
$arr[13,9] = 0.180788118796908
$arr[13,10] = 0.6814665852696464
$arr[14,0] = 0.1303019099362501
$arr[14,1] = 0.07914451678415446
$arr[14,2] = 0.1172338566736322
$arr[14,3] = 2.578061687008969
$arr[14,4] = 1.890997042267173
$arr[14,5] = 1.630106552572713
$arr[14,6] = 1.843343531565665
$arr[14,7] = 0.1847943010804869
$arr[14,8] = 2.452569486421339
$arr[14,9] = 0.1813650096831836
$arr[14,10] = 0.661709091374604
$arr[15,0] = 0.129428948734045
$arr[15,1] = 0.07840071547352778
$arr[15,2] = 0.1172269582962109
$arr[15,3] = 2.582804915389943
$arr[15,4] = 1.895210122912701
$arr[15,5] = 1.63487206609615
$arr[15,6] = 1.846972041075183
$arr[15,7] = 0.1853215057506468
$arr[15,8] = 2.387262187453302
$arr[15,9] = 0.1817324394044348
$arr[15,10] = 0.6496281691038206
$arr[16,0] = 0.1289298631494233
$arr[16,1] = 0.07797368869545096
$arr[16,2] = 0.1172261678174547
$arr[16,3] = 2.58568102593479
$arr[16,4] = 1.897761356296343
$arr[16,5] = 1.637696428545624
$arr[16,6] = 1.849166637409908
$arr[16,7] = 0.1856310495380988
$arr[16,8] = 2.349733040187573
$arr[16,9] = 0.1819487515963019
$arr[16,10] = 0.6426942276156495
$arr[17,0] = 0.128761401174728
$arr[17,1] = 0.07782924115385015
$arr[17,2] = 0.1172264464385577
$arr[17,3] = 2.586680223084699
$arr[17,4] = 1.898647122564029
$arr[17,5] = 1.638667022755683
$arr[17,6] = 1.849928156552451
$arr[17,7] = 0.1857369397330899
$arr[17,8] = 2.337032185652845
$arr[17,9] = 0.1820228464139433
$arr[17,10] = 0.6403490469569491
$arr[18,0] = 0.1295215648569013
$arr[18,1] = 0.07847981307769203
$arr[18,2] = 0.1172273639366512
$arr[18,3] = 2.582284677513726
$arr[18,4] = 1.894748383375315
$arr[18,5] = 1.634356140006446
$arr[18,6] = 1.846574644343626
$arr[18,7] = 0.1852647309648248
$arr[18,8] = 2.394210766013998
$arr[18,9] = 0.1816928109781202
$arr[18,10] = 0.6509126865144381
$arr[19,0] = 0.1321115750268405
$arr[19,1] = 0.08067498508893323
$arr[19,2] = 0.1172683668341552
$arr[19,3] = 2.569284573296841
$arr[19,4] = 1.883179605772042
$arr[19,5] = 1.620868525792517
$arr[19,6] = 1.836593266862344
$arr[19,7] = 0.1837535377694692
$arr[19,8] = 2.58672410958485
$arr[19,9] = 0.1806433539903622
$arr[19,10] = 0.6865797537855372
$arr[20,0] = 0.1338308706836102
$arr[20,1] = 0.08211644510433302
$arr[20,2] = 0.1173233329231191
$arr[20,3] = 2.562108209577573
$arr[20,4] = 1.876762514769453
$arr[20,5] = 1.612796948254513
$arr[20,6] = 1.831029892826116
$arr[20,7] = 0.1828220042852706
$arr[20,8] = 2.712833322139545
$arr[20,9] = 0.1800016676120499
$arr[20,10] = 0.7100186844143792
$arr[21,0] = 0.1329108333813735
$arr[21,1] = 0.08134650062642379
$arr[21,2] = 0.1172914251012855
$arr[21,3] = 2.56581746530378
$arr[21,4] = 1.880082581885205
$arr[21,5] = 1.617036834972495
$arr[21,6] = 1.833911225052994
$arr[21,7] = 0.1833140496214973
$arr[21,8] = 2.645500402920732
$arr[21,9] = 0.180340109656747
$arr[21,10] = 0.6974972454319186
$arr[22,0] = 0.1294796844104269
$arr[22,1] = 0.07844405119799092
$arr[22,2] = 0.1172271706540506
$arr[22,3] = 2.582519412429633
$arr[22,4] = 1.894956733785293
$arr[22,5] = 1.634589126648621
$arr[22,6] = 1.846753969390733
$arr[22,7] = 0.1852903787630922
$arr[22,8] = 2.39106926092137
$arr[22,9] = 0.1817107111893179
$arr[22,10] = 0.6503319202168711
$arr[23,0] = 0.1258942136934849
$arr[23,1] = 0.07534354858618997
$arr[23,2] = 0.1172794090979128
$arr[23,3] = 2.606189274891094
$arr[23,4] = 1.915891412555041
$arr[23,5] = 1.656708762562317
$arr[23,6] = 1.864718444226128
$arr[23,7] = 0.1876634437936033
$arr[23,8] = 2.117925375860807
$arr[23,9] = 0.1833794562699893
$arr[23,10] = 0.6000191773612755

$ws.Range("C2:M25").Value = $arr
